$d = $word.ActiveDocument

# Locate the paragraph that holds the "LOM3246: ... (Requisito)" line.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "*LOM3246*Requisito*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # The three paragraphs right after it are expected to be:
    #   1) a blank paragraph
    #   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   3) "© 2020 . Contact: ... Creative Commons Attribution"
    # Those are the footer/navigation boilerplate paragraphs that the
    # site rebuild dropped from the page. Remove them (but leave the
    # following blank paragraph / page-break paragraph untouched).
    $p1 = $targetIndex + 1
    $p2 = $targetIndex + 2
    $p3 = $targetIndex + 3

    if ($p3 -le $d.Paragraphs.Count) {
        $t1 = $d.Paragraphs.Item($p1).Range.Text
        $t2 = $d.Paragraphs.Item($p2).Range.Text
        $t3 = $d.Paragraphs.Item($p3).Range.Text

        $isBlank = ($t1.Trim() -eq "")
        $isJupiter = ($t2 -like "*Jupiter*")
        $isCopyright = ($t3 -like "*Contact:*")

        if ($isBlank -and $isJupiter -and $isCopyright) {
            $start = $d.Paragraphs.Item($targetIndex).Range.End
            $end = $d.Paragraphs.Item($p3).Range.End
            $r = $d.Range($start, $end)
            $r.Delete()
        }
    }
}
